$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("hpi")
$ws.Range("B3").Value = "Skin changes are characteristic of CREST syndrome and less pronounced in Scleroderma."
$ws.Range("C3").Value = "Internal organ involvement (such as renal or cardiac)"
$ws.Range("D3").Value = "Internal organ involvement is more common in Scleroderma compared to CREST syndrome."
$ws.Range("C4").Value = "Sclerodactyly (thickening and tightening of the skin on fingers)"
$ws.Range("D4").Value = "Sclerodactyly is a hallmark of Scleroderma and not typically seen in CREST syndrome."
$ws.Range("B5").Value = "Gastroesophageal reflux is more prevalent in patients with Type 2 Achalasia, which is associated with CREST syndrome."
$ws.Range("C5").Value = "Digital ulcers"
$ws.Range("D5").Value = "Digital ulcers are more frequently associated with Scleroderma than with CREST syndrome."
$ws.Range("B6").Value = "Pulmonary hypertension is a common complication in CREST syndrome, distinguishing it from Scleroderma."
$ws.Range("C6").Value = "Gastrointestinal motility disorders (not specifically achalasia)"
$ws.Range("D6").Value = "Gastrointestinal motility disorders are more generalized in Scleroderma compared to the specific achalasia in CREST syndrome."
$ws = $wb.Worksheets.Item("hist")
$ws.Range("A2").Value = "History of esophageal dilation procedures"
$ws.Range("B2").Value = "Patients with CREST syndrome and Type 2 Achalasia often undergo esophageal dilation due to dysphagia, which is less common in Scleroderma."
$ws.Range("C2").Value = "History of diffuse skin involvement"
$ws.Range("D2").Value = "Diffuse skin involvement is a classic feature of Scleroderma, distinguishing it from the localized skin changes in CREST syndrome."
$ws.Range("B3").Value = "Raynaud's phenomenon is a hallmark of CREST syndrome and is less frequently reported in Scleroderma."
$ws.Range("A4").Value = "History of skin thickening localized to fingers"
$ws.Range("B4").Value = "Localized skin thickening is more characteristic of CREST syndrome compared to the diffuse skin involvement seen in Scleroderma."
$ws.Range("D4").Value = "Renal crisis is a serious complication that occurs more frequently in Scleroderma compared to CREST syndrome."
$ws.Range("B5").Value = "GERD is more prevalent in patients with Type 2 Achalasia, leading to the use of proton pump inhibitors, which is less common in Scleroderma."
$ws.Range("D5").Value = "Immunosuppressive therapy is more commonly required in Scleroderma due to its systemic nature, unlike in CREST syndrome."
$ws.Range("B6").Value = "Pulmonary hypertension is a common complication in CREST syndrome, while it is less frequently noted in Scleroderma."
$ws.Range("C6").Value = "History of gastrointestinal motility disorders"
$ws.Range("D6").Value = "While gastrointestinal issues can occur in both, motility disorders are more prevalent in Scleroderma, distinguishing it from Type 2 Achalasia."
$ws = $wb.Worksheets.Item("soc")
$ws.Range("C2").Value = "Family history of Scleroderma"
$ws.Range("D2").Value = "A direct family history of Scleroderma is a strong indicator for the diagnosis, as it has a genetic component."
$ws.Range("A3").Value = "Occupational exposure to silica or other environmental toxins"
$ws.Range("C3").Value = "History of Raynaud's phenomenon"
$ws.Range("D3").Value = "Raynaud's phenomenon is more commonly associated with Scleroderma, making it a distinguishing feature."
$ws.Range("A4").Value = "History of esophageal symptoms or dysphagia"
$ws.Range("B4").Value = "Dysphagia is more prevalent in CREST syndrome due to Type 2 Achalasia, distinguishing it from Scleroderma."
$ws.Range("C4").Value = "History of pulmonary fibrosis or interstitial lung disease"
$ws.Range("D4").Value = "These respiratory complications are more frequently seen in Scleroderma, indicating a different disease process."
$ws.Range("B5").Value = "Increased risk of vascular issues is more characteristic of CREST syndrome, which can lead to related complications."
$ws.Range("C5").Value = "History of gastrointestinal motility disorders without esophageal involvement"
$ws.Range("D5").Value = "Gastrointestinal motility issues are common in Scleroderma but may not specifically involve the esophagus as in Type 2 Achalasia."
$ws.Range("A6").Value = "History of limited skin involvement"
$ws.Range("B6").Value = "Limited skin involvement is more typical in CREST syndrome compared to the diffuse skin changes seen in Scleroderma."
$ws.Range("C6").Value = "History of joint pain or arthritis"
$ws.Range("D6").Value = "Joint pain is more prevalent in Scleroderma, which can help differentiate it from CREST syndrome."
$ws = $wb.Worksheets.Item("obj")
$ws.Range("B2").Value = "This finding is more characteristic of CREST syndrome, which presents with localized sclerodactyly."
$ws.Range("C2").Value = "Skin thickening over the trunk and proximal extremities"
$ws.Range("D2").Value = "This finding is indicative of diffuse scleroderma, which is more extensive than the limited skin involvement in CREST syndrome."
$ws.Range("A3").Value = "Raynaud's phenomenon"
$ws.Range("B3").Value = "Raynaud's is more prevalent in CREST syndrome and can be associated with esophageal motility issues."
$ws.Range("C3").Value = "Pulmonary fibrosis"
$ws.Range("D3").Value = "Pulmonary fibrosis is a common complication of scleroderma and less frequently seen in CREST syndrome."
$ws.Range("A4").Value = "Dysphagia (difficulty swallowing)"
$ws.Range("B4").Value = "Dysphagia is common in Type 2 Achalasia due to esophageal motility dysfunction, which is a feature of CREST syndrome."
$ws.Range("C4").Value = "Gastrointestinal motility issues beyond esophagus"
$ws.Range("D4").Value = "Scleroderma can cause widespread gastrointestinal involvement, while CREST syndrome is more localized."
$ws.Range("B5").Value = "Telangiectasia is a hallmark of CREST syndrome and less commonly seen in Scleroderma."
$ws.Range("D5").Value = "Joint contractures are more common in scleroderma due to skin and soft tissue involvement."
$ws.Range("A6").Value = "Calcinosis cutis (calcium deposits in the skin)"
$ws.Range("B6").Value = "Calcinosis is more frequently associated with CREST syndrome than with diffuse scleroderma."
$ws.Range("C6").Value = "Sclerodactyly (thickened skin on fingers)"
$ws.Range("D6").Value = "Sclerodactyly is a classic sign of scleroderma and is more pronounced than in CREST syndrome."
$ws = $wb.Worksheets.Item("test")
$ws.Range("B3").Value = "This imaging finding is more common in patients with Type 2 Achalasia, indicating esophageal motility issues."
$ws.Range("D3").Value = "These antibodies are more commonly associated with diffuse Scleroderma rather than CREST syndrome."
$ws.Range("B4").Value = "These antibodies are more frequently associated with CREST syndrome than with Scleroderma."
$ws.Range("D4").Value = "This finding is more indicative of Scleroderma, which often leads to pulmonary complications."
$ws.Range("B5").Value = "While both conditions can show lung involvement, the pattern is more pronounced in CREST syndrome due to fibrosis."
$ws.Range("C5").Value = "Echocardiogram showing right heart failure or pulmonary hypertension"
$ws.Range("D5").Value = "These cardiac complications are more frequently seen in Scleroderma than in CREST syndrome."
$ws.Range("C6").Value = "Laboratory tests showing renal involvement (elevated creatinine)"
$ws.Range("D6").Value = "Renal crisis is a significant complication of Scleroderma, distinguishing it from CREST syndrome."
